$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2843.4
$ws.Range("I33").Value = 100.181816
$ws.Range("J33").Value = 10387.25
$ws.Range("K33").Value = 100.181816
$ws.Range("L33").Value = 10387.25
$ws.Range("M33").Value = 128.818184
$ws.Range("N33").Value = -10845.25
$ws.Range("H34").Value = 2833.3333
$ws.Range("I34").Value = 2833.3333
$ws.Range("K34").Value = 2833.3333
$ws.Range("M34").Value = -2630.3333
$ws.Range("H36").Value = 2833.3333
$ws.Range("I36").Value = 2833.3333
$ws.Range("K36").Value = 2833.3333
$ws.Range("M36").Value = -2118.3333
$ws.Range("H53").Value = 2855.6667
$ws.Range("J53").Value = 5000
$ws.Range("L53").Value = 5000
$ws.Range("N53").Value = -6274
$ws.Range("H70").Value = 6669400
$ws.Range("J70").Value = 7145610.5
$ws.Range("L70").Value = 21436831.5
$ws.Range("N70").Value = -21437371.5
$ws.Range("H73").Value = 6669400
$ws.Range("J73").Value = 7145610.5
$ws.Range("L73").Value = 21436831.5
$ws.Range("N73").Value = -21438703.5
$ws.Range("H132").Value = 1815.0244
$ws.Range("I132").Value = 1382.4872
$ws.Range("J132").Value = 10249.5
$ws.Range("K132").Value = 4147.461600000001
$ws.Range("L132").Value = 30748.5
$ws.Range("M132").Value = -1617.461600000001
$ws.Range("N132").Value = -35808.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9798.75
$ws.Range("I2").Value = 1598.6364
$ws.Range("K2").Value = 1598.6364
$ws.Range("M2").Value = -1485.6364
$ws.Range("H32").Value = 1840.3462
$ws.Range("I32").Value = 1889.5416
$ws.Range("K32").Value = 1889.5416
$ws.Range("M32").Value = -1602.5416
$ws.Range("H45").Value = 66671164
$ws.Range("I45").Value = 111114504
$ws.Range("J45").Value = 6154.6665
$ws.Range("K45").Value = 111114504
$ws.Range("L45").Value = 6154.6665
$ws.Range("M45").Value = -111114127
$ws.Range("N45").Value = -6908.6665
$ws.Range("H63").Value = 4305.7407
$ws.Range("I63").Value = 2448.3572
$ws.Range("J63").Value = 6306
$ws.Range("K63").Value = 2448.3572
$ws.Range("L63").Value = 6306
$ws.Range("M63").Value = -1762.3572
$ws.Range("N63").Value = -7678
$ws.Range("H66").Value = 4305.7407
$ws.Range("I66").Value = 2448.3572
$ws.Range("J66").Value = 6306
$ws.Range("K66").Value = 12241.786
$ws.Range("L66").Value = 31530
$ws.Range("M66").Value = -8809.786
$ws.Range("N66").Value = -38394
$ws.Range("H68").Value = 45000
$ws.Range("J68").Value = 45000
$ws.Range("L68").Value = 45000
$ws.Range("H71").Value = 45000
$ws.Range("J71").Value = 45000
$ws.Range("L71").Value = 135000
$ws.Range("H97").Value = 5184.5
$ws.Range("I97").Value = 4269.3335
$ws.Range("J97").Value = 7930
$ws.Range("K97").Value = 4269.3335
$ws.Range("L97").Value = 7930
$ws.Range("M97").Value = -3773.3335
$ws.Range("N97").Value = -8922
$ws.Range("H102").Value = 3443.2856
$ws.Range("I102").Value = 3350.5
$ws.Range("K102").Value = 3350.5
$ws.Range("M102").Value = -1728.5
$ws.Range("H110").Value = 3423.8
$ws.Range("I110").Value = 2206.5
$ws.Range("K110").Value = 2206.5
$ws.Range("M110").Value = -161.5
$ws.Range("H113").Value = 63157
$ws.Range("J113").Value = 63157
$ws.Range("L113").Value = 63157
$ws.Range("H116").Value = 9798.75
$ws.Range("I116").Value = 1598.6364
$ws.Range("K116").Value = 1598.6364
$ws.Range("M116").Value = 695.3635999999999
$ws.Range("H132").Value = 4684.0444
$ws.Range("I132").Value = 3788.3845
$ws.Range("J132").Value = 10505.833
$ws.Range("K132").Value = 11365.1535
$ws.Range("L132").Value = 31517.499
$ws.Range("M132").Value = -8835.1535
$ws.Range("N132").Value = -36577.499
$ws.Range("N68").Value = -46622
$ws.Range("N71").Value = -143112
$ws.Range("N113").Value = -71835

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9798.75
$ws.Range("I3").Value = 1598.6364
$ws.Range("K3").Value = 1598.6364
$ws.Range("M3").Value = -1484.6364
$ws.Range("H11").Value = 1183.5714
$ws.Range("I11").Value = 1502.375
$ws.Range("J11").Value = 758.5
$ws.Range("K11").Value = 1502.375
$ws.Range("L11").Value = 758.5
$ws.Range("M11").Value = -1362.375
$ws.Range("N11").Value = -1038.5
$ws.Range("H105").Value = 18321.688
$ws.Range("I105").Value = 22994.9
$ws.Range("J105").Value = 10533
$ws.Range("K105").Value = 22994.9
$ws.Range("L105").Value = 10533
$ws.Range("M105").Value = -21247.9
$ws.Range("N105").Value = -14027
$ws.Range("H107").Value = 2940.28
$ws.Range("I107").Value = 2544.75
$ws.Range("K107").Value = 2544.75
$ws.Range("M107").Value = -624.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24251.037
$ws.Range("I31").Value = 3039.8823
$ws.Range("J31").Value = 34267.418
$ws.Range("K31").Value = 3039.8823
$ws.Range("L31").Value = 34267.418
$ws.Range("M31").Value = -2744.8823
$ws.Range("N31").Value = -34857.418
$ws.Range("H34").Value = 24251.037
$ws.Range("I34").Value = 3039.8823
$ws.Range("J34").Value = 34267.418
$ws.Range("K34").Value = 3039.8823
$ws.Range("L34").Value = 34267.418
$ws.Range("M34").Value = -2837.8823
$ws.Range("N34").Value = -34671.418
$ws.Range("H132").Value = 3064.2058
$ws.Range("I132").Value = 2412.724
$ws.Range("J132").Value = 6842.8
$ws.Range("K132").Value = 7238.172
$ws.Range("L132").Value = 20528.4
$ws.Range("M132").Value = -4708.172
$ws.Range("N132").Value = -25588.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 7221.6665
$ws.Range("I69").Value = 3246
$ws.Range("J69").Value = 9209.5
$ws.Range("K69").Value = 9738
$ws.Range("L69").Value = 27628.5
$ws.Range("M69").Value = -8927
$ws.Range("N69").Value = -29250.5
$ws.Range("H72").Value = 7221.6665
$ws.Range("I72").Value = 3246
$ws.Range("J72").Value = 9209.5
$ws.Range("K72").Value = 29214
$ws.Range("L72").Value = 82885.5
$ws.Range("M72").Value = -25158
$ws.Range("N72").Value = -90997.5
$ws.Range("H127").Value = 2152.4
$ws.Range("J127").Value = 2152.4
$ws.Range("L127").Value = 6457.200000000001
$ws.Range("N127").Value = -16377.2
$ws.Range("H138").Value = 8099.75
$ws.Range("I138").Value = 2466.3333
$ws.Range("K138").Value = 7398.999899999999
$ws.Range("M138").Value = -2258.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1554.5834
$ws.Range("I107").Value = 648.2941
$ws.Range("J107").Value = 3755.5715
$ws.Range("K107").Value = 648.2941
$ws.Range("L107").Value = 3755.5715
$ws.Range("M107").Value = 1271.7059
$ws.Range("N107").Value = -7595.5715
$ws.Range("H114").Value = 45000
$ws.Range("I114").Value = 45000
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 45000
$ws.Range("L114").Value = 0
$ws.Range("M114").ClearContents()
$ws.Range("N114").Value = -40661
$ws.Range("H126").Value = 7429.5
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 11501.625
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 34504.875
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -39444.875
$ws.Range("H132").Value = 3153.1853
$ws.Range("I132").Value = 3161.0386
$ws.Range("K132").Value = 9483.1158
$ws.Range("M132").Value = -6953.1158
$ws.Range("H136").Value = 69319
$ws.Range("J136").Value = 69319
$ws.Range("L136").Value = 207957
$ws.Range("N136").Value = -213057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 100000
$ws.Range("J36").Value = 100000
$ws.Range("L36").Value = 100000
$ws.Range("H105").Value = 28000
$ws.Range("J105").Value = 28000
$ws.Range("L105").Value = 28000
$ws.Range("N105").Value = -34988
$ws.Range("H133").Value = 96867.625
$ws.Range("J133").Value = 96867.625
$ws.Range("L133").Value = 96867.625
$ws.Range("N133").Value = -101927.625
$ws.Range("H136").Value = 6586.514
$ws.Range("J136").Value = 9940.6
$ws.Range("L136").Value = 29821.8
$ws.Range("N136").Value = -34921.8
$ws.Range("N36").Value = -101124

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 49950
$ws.Range("J46").Value = 49950
$ws.Range("L46").Value = 49950
$ws.Range("N46").Value = -50412
$ws.Range("H74").Value = 9000.333
$ws.Range("J74").Value = 9000.333
$ws.Range("L74").Value = 9000.333
$ws.Range("N74").Value = -10872.333
$ws.Range("H77").Value = 9000.333
$ws.Range("J77").Value = 9000.333
$ws.Range("L77").Value = 27000.999
$ws.Range("N77").Value = -36360.999
$ws.Range("H107").Value = 834.4706
$ws.Range("I107").Value = 978.6
$ws.Range("K107").Value = 2935.8
$ws.Range("M107").Value = -1015.8
$ws.Range("H117").Value = 68421
$ws.Range("J117").Value = 68421
$ws.Range("L117").Value = 68421
$ws.Range("N117").Value = -77599
$ws.Range("H126").Value = 4270.591
$ws.Range("I126").Value = 3130.6667
$ws.Range("J126").Value = 5059.769
$ws.Range("K126").Value = 9392.000100000001
$ws.Range("L126").Value = 15179.307
$ws.Range("M126").Value = -6922.000100000001
$ws.Range("N126").Value = -20119.307
$ws.Range("H134").Value = 49950
$ws.Range("J134").Value = 49950
$ws.Range("L134").Value = 149850
$ws.Range("N134").Value = -154920
$ws.Range("H136").Value = 5795.2085
$ws.Range("I136").Value = 4934.4287
$ws.Range("J136").Value = 7000.3
$ws.Range("K136").Value = 14803.2861
$ws.Range("L136").Value = 21000.9
$ws.Range("M136").Value = -12253.2861
$ws.Range("N136").Value = -26100.9
